$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 835
$ws1.Range("F16").Value = 7468
$ws1.Range("F28").Value = 6169
$ws1.Range("F33").Value = 446
$ws1.Range("F34").Value = 6404
$ws1.Range("F46").Value = 432

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F17").Value = 7468
$ws4.Range("F29").Value = 6169
$ws4.Range("F35").Value = 446
$ws4.Range("F36").Value = 6404
$ws4.Range("F46").Value = 432
